# Add a "学费" (tuition fee) column (H) to the sheet, with a header cell
# styled like the other header cells, and per-row tuition values that
# follow the existing 专业/学习方式 (program) grouping:
#   rows 2-62   -> 32000
#   rows 63-103 -> 15000
#   rows 104-239-> 28000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell H1 -------------------------------------------------
$ws.Cells.Item(1, 8).Value = "学费"

$g1 = $ws.Range("G1")
$h1 = $ws.Range("H1")
$g1.Copy()
$h1.PasteSpecial(-4122)
$h1.Font.Name = "宋体"

# --- Data cells H2:H239 ----------------------------------------------
for ($r = 2; $r -le 239; $r++) {
    if ($r -le 62) {
        $v = 32000
    } elseif ($r -le 103) {
        $v = 15000
    } else {
        $v = 28000
    }
    $ws.Cells.Item($r, 8).Value = $v
}

# Copy formatting (border/alignment/font) from an existing data cell so
# the new column matches the rest of the table visually.
$g2 = $ws.Range("G2")
$hData = $ws.Range("H2:H239")
$g2.Copy()
$hData.PasteSpecial(-4122)

# --- Selection / scroll position matches the author's saved view -----
$ws.Application.ActiveWindow.ScrollRow = 222
$ws.Range("H104:H239").Select()
